$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Rohit Menon"
$ws.Range("B2").Value = "Posts"
$ws.Range("C2").Value = "Rejected"
$ws.Range("D2").Value = "John Smith"
$ws.Range("E2").Value = "Comments"

# Row 3
$ws.Range("A3").Value = "Rohit Menon"
$ws.Range("B3").Value = "Comments"
$ws.Range("C3").Value = "Post Error"
$ws.Range("D3").Value = "John Smith"
$ws.Range("E3").Value = "Posts"

# Row 4
$ws.Range("A4").Value = "Rohit Menon"
$ws.Range("B4").Value = "Replies"
$ws.Range("C4").Value = "Approved Scheduled"
$ws.Range("D4").Value = "John Smith"
$ws.Range("E4").Value = "Replies"

# Row 5
$ws.Range("A5").Value = "Rohit Menon"
$ws.Range("B5").Value = "All"
$ws.Range("C5").Value = "Pending Approval"
$ws.Range("D5").Value = "John Smith"
$ws.Range("E5").Value = "All"

# Column C width adjustment (bestFit width changed from 18.43 to 19.71
# to accommodate the new, longer values such as "Approved Scheduled")
$ws.Columns.Item(3).ColumnWidth = 19.7109375

# Update selection to H1
$ws.Range("H1").Select()
